# Offer import now support investor offers
#
# Inserts two new columns after "First Name *" (A) to capture whether the
# offeree is a Founder / Employee / Investor, and (for the data row that
# already carried an email address) turns that email into a clickable
# mailto: hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before column B. Everything that used to live in
# B:H shifts right to D:J (and the two "extra"/reserved trailing columns
# shift from G:H to I:J, etc.)
$ws.Columns("B:C").Insert()

# Fill in the new "Founder/Employee/Investor" column for the existing rows.
# (Write the data cells before the header cells so the shared-string table
# order matches: "Founder" then "Employee" get interned first, followed by
# the header text itself.)
$ws.Range("B2").Value = "Founder"
$ws.Range("B3").Value = "Employee"
$ws.Range("B4").Value = "Employee"
$ws.Range("B5").Value = "Employee"

# Header row for the two new columns.
$ws.Range("B1").Value = "Founder/Employee/Investor *"
$ws.Range("C1").Value = "Investor "

# The third offer row's email address (now in column D after the insert)
# becomes a live mailto: link.
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:emp3@myfirm.com")

# Match the slightly narrower width used for the two new columns.
$ws.Columns("B:C").ColumnWidth = 27.42

# Leave the selection where the author left it.
$ws.Range("C4").Select() | Out-Null
